$d = $word.ActiveDocument

# 1) Remove the "Meta description: ..." paragraph that follows the
#    "Play Bill & Ted's Excellent Adventure Slot for Free" Heading1 paragraph.
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# 2) Insert a new bold paragraph ("Play Bill & Ted's Excellent Adventure
#    Slot for Free") right before the final paragraph (the one that used to
#    read "Create a feature image for ..."), matching the structure of the
#    paragraph removed in step 1 (leading empty run + bold run).
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertionPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)
$xmlNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$newParaXml = "<w:p $xmlNs><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Bill &amp; Ted's Excellent Adventure Slot for Free</w:t></w:r></w:p><w:p $xmlNs/>"
$insertionPoint.InsertXML($newParaXml)

# InsertXML leaves behind one extra empty paragraph mark (because the final
# mark of the inserted fragment merges with the following paragraph); remove it.
# (Paragraph.Range.Text of an empty paragraph is just the paragraph mark
# character, so compare the trimmed text rather than an empty string.)
$emptyParaIndex = $d.Paragraphs.Count - 1
$emptyPara = $d.Paragraphs.Item($emptyParaIndex)
if ($emptyPara.Range.Text.Trim() -eq "") {
    $emptyPara.Range.Delete()
}

# 3) Replace the text of the final paragraph (still italic) with the new
#    "Read our review ..." text that used to live in the Meta description.
#    Use a plain Range.Text assignment (rather than Find/Replace) so that no
#    autocorrect/smart-quote substitution touches the apostrophes.
$finalPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$finalStart = $finalPara.Range.Start
$finalEnd = $finalPara.Range.End - 1   # exclude the trailing paragraph mark
$finalTextRange = $d.Range($finalStart, $finalEnd)
$finalTextRange.Text = "Read our review of Bill & Ted's Excellent Adventure slot machine, play for free, and enjoy fun gameplay features and real cast images."
